$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting (style s="2") from the last existing row (A229)
# down onto the four new date cells so the new rows match the existing look.
$ws.Range("A229").Copy()
$ws.Range("A230:A233").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 230
$ws.Range("A230").Value = 44304
$ws.Range("B230").Value = 0
$ws.Range("C230").Value = 2
$ws.Range("D230").Value = 166.8056713928273

# Row 231
$ws.Range("A231").Value = 44305
$ws.Range("B231").Value = 0
$ws.Range("C231").Value = 0
$ws.Range("D231").Value = 0

# Row 232
$ws.Range("A232").Value = 44306
$ws.Range("B232").Value = 0
$ws.Range("C232").Value = 0
$ws.Range("D232").Value = 0

# Row 233
$ws.Range("A233").Value = 44307
$ws.Range("B233").Value = 0
$ws.Range("C233").Value = 0
$ws.Range("D233").Value = 0
